$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.43820000000003
$ws.Range("D5").Value = -8.522099999999995
$ws.Range("D9").Value = -8.722900000000005
$ws.Range("D11").Value = -8.316600000000003
$ws.Range("A21").Value = -21.32740000000001
$ws.Range("D21").Value = -8.052900000000005
$ws.Range("A23").Value = -21.41640000000003
$ws.Range("A25").Value = -22.35440000000003
